$wb = $excel.ActiveWorkbook

# "Table - 2.1" now has answers filled in for the average funding
# amount questions (rows 5-8, column C).
$ws = $wb.Worksheets.Item("Table - 2.1")
$ws.Range("C5").Value = 11748949.1
$ws.Range("C6").Value = 958694.5
$ws.Range("C7").Value = 719818
$ws.Range("C8").Value = 73308593

# Work moved on to "Table - 2.1": it becomes the active/selected sheet
# (this also drops tabSelected from "Table -1.1"), with C19 selected.
$ws.Activate()
$ws.Range("C19").Select()

# Widen the app window a bit, as recorded in the workbook view.
$excel.ActiveWindow.Width = 28800
